$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.715.32'
$ws.Range("D3").Value = '1.601.48'
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.55'
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.512'
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.246'
$ws.Range("E9").Value = '  +0.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.73'
$ws.Range("E10").Value = '  +0.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0841'
$ws.Range("E11").Value = '  +0.56%  '
$ws.Range("D12").Value = '1.825.46'
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").Value = '1.605.24'
$ws.Range("E13").Value = '  +1.76%  '
$ws.Range("E14").Value = '  +0.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.524'
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.18'
$ws.Range("E16").Value = '  +0.27%  '
$ws.Range("D17").Value = '26.683.10'
$ws.Range("D18").Value = '0.0₃0748'
$ws.Range("E18").Value = '  +1.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.25'
$ws.Range("E19").Value = '  +2.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '210.22'
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("E21").Value = '  +0.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.31'
$ws.Range("E22").Value = '  +0.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.30'
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.97'
$ws.Range("E24").Value = '  +0.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.02'
$ws.Range("E25").Value = '  -1.92%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  -1.34%  '
$ws.Range("E28").Value = '  -1.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.43'
$ws.Range("E29").Value = '  +0.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0515'
$ws.Range("E30").Value = '  +1.19%  '
$ws.Range("E31").Value = '  -0.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.28'
$ws.Range("E32").Value = '  +1.29%  '
$ws.Range("D34").Value = '1.293.06'
$ws.Range("E34").Value = '  +0.62%  '
$ws.Range("E35").Value = '  +0.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.611'
$ws.Range("E36").Value = '  -3.85%  '
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.14'
$ws.Range("E38").Value = '  +21.47%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0170'
$ws.Range("E39").Value = '  -0.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.825'
$ws.Range("E40").Value = '  -2.57%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.24'
$ws.Range("E41").Value = '  +1.49%  '
$ws.Range("E42").Value = '  -1.77%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.781'
$ws.Range("E43").Value = '  -1.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.20'
$ws.Range("E44").Value = '  -1.62%  '
$ws.Range("D45").Value = '1.736.44'
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.31'
$ws.Range("E46").Value = '  +1.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.56'
$ws.Range("E47").Value = '  -3.06%  '
$ws.Range("D48").Value = '0.0₆0105'
$ws.Range("E48").Value = '  -0.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.102'
$ws.Range("E49").Value = '  -1.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0516'
$ws.Range("E50").Value = '  +1.53%  '
$ws.Range("E51").Value = '  +0.24%  '
